$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 26 (existing rows 26:33 shift down to 27:34),
# and populate it with this week's new price record for Espárragos.
$ws.Rows("26:26").Insert()

$ws.Range("A26").Value = 7
$ws.Range("B26").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C26").Value = "Ñuble"
$ws.Range("D26").Value = 44876
$ws.Range("E26").Value = 16
$ws.Range("F26").Value = 300000000
$ws.Range("G26").Value = "Espárragos"
$ws.Range("H26").Value = "Sin especificar"
$ws.Range("I26").Value = "Primera"
$ws.Range("J26").Value = 1000
$ws.Range("K26").Value = 1000
$ws.Range("L26").Value = 1100
$ws.Range("M26").Value = 1050
$ws.Range("N26").Value = "$/kilo"
$ws.Range("O26").Value = "Región de Ñuble"
$ws.Range("P26").Value = 1050
$ws.Range("Q26").Value = 1
$ws.Range("R26").Value = "Hortaliza"
